# "Updated symbol list" refresh of the crypto price table on Sheet1.
#
# Two kinds of edits happen here:
#  1. Plain price (column D) refreshes for rows that keep the same coin.
#  2. Rows 10-18: the "One" coin moved up into the top-10 list, so every
#     coin from row 10 down to row 18 shifted down by one slot (each row's
#     name/link/price/volume now shows what used to be one row below it),
#     and "One" was inserted at row 10 with a new price.
#
# Column D values are stored as plain text in the workbook (e.g. "229.90",
# keeping trailing zeros), so a leading apostrophe is used to force text
# entry and stop them from being auto-converted to numbers (which would
# silently drop significant trailing zeros, e.g. "3.260" -> 3.26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- simple price-only refreshes -----------------------------------------
$ws.Range("D2").Value  = "'229.90"
$ws.Range("D3").Value  = "'22.35"
$ws.Range("D4").Value  = "'5.268"
$ws.Range("D7").Value  = "'6.480"
$ws.Range("D8").Value  = "'1.057"
$ws.Range("D9").Value  = "'0.7817"

# --- rows 10-18: coin list shifted down by one (One entered at row 10) ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.0005900"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1381"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07393"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03151"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02972"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09271"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001663"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "'3.260"
$ws.Range("E17").Value = "16MCDexMCB"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04774"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- more simple price-only refreshes -------------------------------------
$ws.Range("D19").Value = "'0.006234"
$ws.Range("D20").Value = "'0.005236"
$ws.Range("D21").Value = "'0.001063"
$ws.Range("D23").Value = "'3.918"

# row 27: price refresh + trailing "Bestin24h" tag removed from volume label
$ws.Range("D27").Value = "'0.0005000"
$ws.Range("E27").Value = "26UpBotsUBXT"

$ws.Range("D40").Value = "'0.04006"

# row 41: price refresh + trailing "Bestin24h" tag added to volume label
$ws.Range("D41").Value = "'0.007036"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

$ws.Range("D42").Value = "'0.003408"
$ws.Range("D44").Value = "'0.009966"
$ws.Range("D48").Value = "'0.04135"
